$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New nomenclature rows appended to the Spanish/English table (A: ES, B: EN)
# Values are written in the order that reproduces the original shared-string
# interning order (35=Execute,36=Ejecución,37=Enlace,38=Link,39=Design,
# 40=Diseño,41=Group,42=grupo).
$ws.Range("B9").Value = "Execute"
$ws.Range("A9").Value = "Ejecución "

$ws.Range("A10").Value = "Enlace"
$ws.Range("B10").Value = "Link"

$ws.Range("B11").Value = "Design"
$ws.Range("A11").Value = "Diseño"

$ws.Range("A12").Value = "Group "
$ws.Range("B12").Value = "grupo "

$ws.Range("A13").Select()
